# Applies the changes described by the commit to abe_replication.xlsx
# Workbook has 4 sheets: "Table 1", "Table 2", "Table 3", "Table 4"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 1: add two new descriptive-stats rows (Age, Gender)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

# Duplicate the formatting of the last existing row (row 5) onto the two new rows
$ws1.Range("A5").Copy($ws1.Range("A6"))
$ws1.Range("A5").Copy($ws1.Range("A7"))

$ws1.Range("A6").Value = "Age"
$ws1.Range("B6").Value = 39.37802291047942
$ws1.Range("C6").Value = 13.41103145474985
$ws1.Range("D6").Value = 17
$ws1.Range("E6").Value = 90

$ws1.Range("A7").Value = "Gender (0: F | 1: M)"
$ws1.Range("B7").Value = 0.5803988120492151
$ws1.Range("C7").Value = 0.4935984192187962
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 1

# ---------------------------------------------------------------------------
# Table 2: update a few validation/calibration numbers for the HB M2 column
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table 2")

$ws2.Range("D3").Value = 0.5600000000000001
$ws2.Range("D4").Value = 0.98
$ws2.Range("D6").Value = 3.06
$ws2.Range("D7").Value = 0.76

# ---------------------------------------------------------------------------
# Table 3: HB M2 now has 3 covariates instead of 1, so 4 new parameter rows
# are inserted (2 for the purchase-rate block, 2 for the dropout-rate block)
# and every posterior-quantile value in the HB M2 columns (E:G) is refreshed.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

$ws3.Range("E1").Value = "HB M2 (with 3 covariates)"

# Insert two rows right after row 5 (Purchase rate log(lambda) - ... block)
$ws3.Range("A6:A7").EntireRow.Insert()
# Insert two rows right after row 9 (Dropout rate log(mu) - ... block, after the first insert)
$ws3.Range("A10:A11").EntireRow.Insert()

# Give the four brand-new label cells the same look (border/bold/alignment) as
# the other row-label cells in column A
$ws3.Range("A4").Copy($ws3.Range("A6"))
$ws3.Range("A4").Copy($ws3.Range("A7"))
$ws3.Range("A4").Copy($ws3.Range("A10"))
$ws3.Range("A4").Copy($ws3.Range("A11"))

# Row 4: Purchase rate log(lambda) - Intercept (HB M2 values refreshed)
$ws3.Range("E4").Value = -3.86
$ws3.Range("F4").Value = -3.64
$ws3.Range("G4").Value = -3.44

# Row 5: relabeled, HB M2 values refreshed
$ws3.Range("A5").Value = "Purchase rate log(λ) - first.sales"
$ws3.Range("E5").Value = 0.07000000000000001
$ws3.Range("F5").Value = 0.21
$ws3.Range("G5").Value = 0.32

# Row 6 (new): Purchase rate log(lambda) - age scaled
$ws3.Range("A6").Value = "Purchase rate log(λ) - age scaled"
$ws3.Range("E6").Value = -0.26
$ws3.Range("F6").Value = -0.11
$ws3.Range("G6").Value = 0.04

# Row 7 (new): Purchase rate log(lambda) - gender binary
$ws3.Range("A7").Value = "Purchase rate log(λ) - gender binary"
$ws3.Range("E7").Value = -0.12
$ws3.Range("F7").Value = 0.08
$ws3.Range("G7").Value = 0.29

# Row 8: Dropout rate log(mu) - Intercept (HB M2 values refreshed)
$ws3.Range("E8").Value = -4.61
$ws3.Range("F8").Value = -3.96
$ws3.Range("G8").Value = -3.53

# Row 9: relabeled, HB M2 values refreshed
$ws3.Range("A9").Value = "Dropout rate log(μ) - first.sales"
$ws3.Range("E9").Value = -0.35
$ws3.Range("F9").Value = 0.06
$ws3.Range("G9").Value = 0.27

# Row 10 (new): Dropout rate log(mu) - age scaled
$ws3.Range("A10").Value = "Dropout rate log(μ) - age scaled"
$ws3.Range("E10").Value = -0.16
$ws3.Range("F10").Value = 0.1
$ws3.Range("G10").Value = 0.29

# Row 11 (new): Dropout rate log(mu) - gender binary
$ws3.Range("A11").Value = "Dropout rate log(μ) - gender binary"
$ws3.Range("E11").Value = -0.14
$ws3.Range("F11").Value = 0.41
$ws3.Range("G11").Value = 1.01

# Row 12 (was row 8): sigma^2_lambda = var[log lambda]
$ws3.Range("F12").Value = 1.4
$ws3.Range("G12").Value = 1.72

# Row 13 (was row 9): sigma^2_mu = var[log mu]
$ws3.Range("E13").Value = -0.21
$ws3.Range("F13").Value = 0.26
$ws3.Range("G13").Value = 0.82

# Row 14 (was row 10): sigma_lambda_mu = cov[log lambda, log mu]
$ws3.Range("E14").Value = 1.11
$ws3.Range("F14").Value = 2.5
$ws3.Range("G14").Value = 4.03

# Row 15 (was row 11): Correlation computed from Gamma0
$ws3.Range("E15").Value = -0.14
$ws3.Range("F15").Value = 0.14
$ws3.Range("G15").Value = 0.35

# Row 16 (was row 12): Marginal log-likelihood
$ws3.Range("F16").Value = -15070

# ---------------------------------------------------------------------------
# Table 4: refresh the per-customer / summary posterior estimates
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table 4")

$ws4.Range("B2").Value = 0.753
$ws4.Range("E2").Value = 0.0163
$ws4.Range("F2").Value = 0.0009
$ws4.Range("G2").Value = 0.0683
$ws4.Range("H2").Value = 1.18
$ws4.Range("I2").Value = 0.429
$ws4.Range("K2").Value = 21.69

$ws4.Range("C3").Value = 0.459
$ws4.Range("D3").Value = 0.988
$ws4.Range("E3").Value = 0.0166
$ws4.Range("F3").Value = 0.0011
$ws4.Range("G3").Value = 0.065
$ws4.Range("H3").Value = 1.16
$ws4.Range("I3").Value = 0.421
$ws4.Range("J3").Value = 0.996
$ws4.Range("K3").Value = 19.93

$ws4.Range("B4").Value = 0.498
$ws4.Range("C4").Value = 0.292
$ws4.Range("D4").Value = 0.758
$ws4.Range("E4").Value = 0.0136
$ws4.Range("G4").Value = 0.0574
$ws4.Range("H4").Value = 1.42
$ws4.Range("I4").Value = 0.494
$ws4.Range("J4").Value = 0.998
$ws4.Range("K4").Value = 15.05

$ws4.Range("B5").Value = 0.515
$ws4.Range("C5").Value = 0.291
$ws4.Range("D5").Value = 0.798
$ws4.Range("E5").Value = 0.0152
$ws4.Range("F5").Value = 0.0008
$ws4.Range("G5").Value = 0.0645
$ws4.Range("H5").Value = 1.26
$ws4.Range("I5").Value = 0.453
$ws4.Range("J5").Value = 0.992
$ws4.Range("K5").Value = 15.03

$ws4.Range("B6").Value = 0.448
$ws4.Range("C6").Value = 0.26
$ws4.Range("D6").Value = 0.6860000000000001
$ws4.Range("E6").Value = 0.0138
$ws4.Range("F6").Value = 0.0007
$ws4.Range("G6").Value = 0.0568
$ws4.Range("H6").Value = 1.4
$ws4.Range("I6").Value = 0.489
$ws4.Range("J6").Value = 0.988
$ws4.Range("K6").Value = 13.36

$ws4.Range("B7").Value = 0.394
$ws4.Range("C7").Value = 0.209
$ws4.Range("D7").Value = 0.638
$ws4.Range("E7").Value = 0.0146
$ws4.Range("G7").Value = 0.0634
$ws4.Range("H7").Value = 1.31
$ws4.Range("I7").Value = 0.467
$ws4.Range("J7").Value = 0.983
$ws4.Range("K7").Value = 11.5

$ws4.Range("D8").Value = 0.618
$ws4.Range("E8").Value = 0.0153
$ws4.Range("F8").Value = 0.0008
$ws4.Range("G8").Value = 0.0655
$ws4.Range("H8").Value = 1.25
$ws4.Range("I8").Value = 0.45
$ws4.Range("J8").Value = 0.974
$ws4.Range("K8").Value = 10.8

$ws4.Range("B9").Value = 0.322
$ws4.Range("C9").Value = 0.167
$ws4.Range("D9").Value = 0.533
$ws4.Range("E9").Value = 0.0133
$ws4.Range("G9").Value = 0.0548
$ws4.Range("H9").Value = 1.45
$ws4.Range("I9").Value = 0.501
$ws4.Range("J9").Value = 0.992
$ws4.Range("K9").Value = 9.710000000000001

$ws4.Range("B10").Value = 0.328
$ws4.Range("C10").Value = 0.171
$ws4.Range("D10").Value = 0.535
$ws4.Range("E10").Value = 0.0165
$ws4.Range("F10").Value = 0.001
$ws4.Range("G10").Value = 0.0654
$ws4.Range("H10").Value = 1.17
$ws4.Range("I10").Value = 0.424
$ws4.Range("J10").Value = 0.992
$ws4.Range("K10").Value = 9.359999999999999

$ws4.Range("B11").Value = 0.303
$ws4.Range("C11").Value = 0.152
$ws4.Range("D11").Value = 0.503
$ws4.Range("E11").Value = 0.0143
$ws4.Range("F11").Value = 0.0007
$ws4.Range("G11").Value = 0.0595
$ws4.Range("H11").Value = 1.35
$ws4.Range("I11").Value = 0.475
$ws4.Range("J11").Value = 0.981
$ws4.Range("K11").Value = 8.890000000000001

$ws4.Range("B13").Value = 0.026
$ws4.Range("C13").Value = 0.001
$ws4.Range("D13").Value = 0.129
$ws4.Range("E13").Value = 0.035
$ws4.Range("F13").Value = 0.0019
$ws4.Range("G13").Value = 1.2343
$ws4.Range("H13").Value = 0.55
$ws4.Range("I13").Value = 0.162
$ws4.Range("J13").Value = 0.254
$ws4.Range("K13").Value = 0.14

$ws4.Range("B14").Value = 0.028
$ws4.Range("C14").Value = 0.001
$ws4.Range("D14").Value = 0.133
$ws4.Range("E14").Value = 0.0352
$ws4.Range("F14").Value = 0.0018
$ws4.Range("G14").Value = 1.3254
$ws4.Range("H14").Value = 0.55
$ws4.Range("I14").Value = 0.16
$ws4.Range("J14").Value = 0.242
$ws4.Range("K14").Value = 0.14

$ws4.Range("B15").Value = 0.492
$ws4.Range("C15").Value = 0.188
$ws4.Range("D15").Value = 0.9379999999999999
$ws4.Range("E15").Value = 0.0396
$ws4.Range("F15").Value = 0.0071
$ws4.Range("G15").Value = 0.2426
$ws4.Range("H15").Value = 0.49
$ws4.Range("I15").Value = 0.128
$ws4.Range("J15").Value = 0.014
$ws4.Range("K15").Value = 0.14

$ws4.Range("B16").Value = 0.029
$ws4.Range("C16").Value = 0.001
$ws4.Range("D16").Value = 0.145
$ws4.Range("E16").Value = 0.0358
$ws4.Range("F16").Value = 0.0019
$ws4.Range("G16").Value = 1.4293
$ws4.Range("H16").Value = 0.54
$ws4.Range("I16").Value = 0.155
$ws4.Range("J16").Value = 0.23
$ws4.Range("K16").Value = 0.14

$ws4.Range("B17").Value = 0.353
$ws4.Range("C17").Value = 0.082
$ws4.Range("D17").Value = 0.84
$ws4.Range("E17").Value = 0.0389
$ws4.Range("F17").Value = 0.0059
$ws4.Range("G17").Value = 0.2982
$ws4.Range("H17").Value = 0.49
$ws4.Range("I17").Value = 0.132
$ws4.Range("J17").Value = 0.019
$ws4.Range("K17").Value = 0.13

$ws4.Range("B18").Value = 0.379
$ws4.Range("C18").Value = 0.07000000000000001
$ws4.Range("D18").Value = 1.009
$ws4.Range("E18").Value = 0.043
$ws4.Range("F18").Value = 0.0083
$ws4.Range("G18").Value = 0.4693
$ws4.Range("H18").Value = 0.45
$ws4.Range("I18").Value = 0.107
$ws4.Range("K18").Value = 0.08

$ws4.Range("B19").Value = 0.421
$ws4.Range("C19").Value = 0.091
$ws4.Range("D19").Value = 1.035
$ws4.Range("E19").Value = 0.0434
$ws4.Range("F19").Value = 0.0098
$ws4.Range("G19").Value = 0.4399
$ws4.Range("H19").Value = 0.44
$ws4.Range("I19").Value = 0.105
$ws4.Range("J19").Value = 0.007
$ws4.Range("K19").Value = 0.06

$ws4.Range("B20").Value = 0.924
$ws4.Range("C20").Value = 0.145
$ws4.Range("D20").Value = 2.688
$ws4.Range("F20").Value = 0.0118
$ws4.Range("G20").Value = 1.0989

$ws4.Range("B21").Value = 0.708
$ws4.Range("C21").Value = 0.223
$ws4.Range("D21").Value = 1.49
$ws4.Range("E21").Value = 0.0411
$ws4.Range("F21").Value = 0.0073
$ws4.Range("G21").Value = 0.3828
$ws4.Range("H21").Value = 0.47
$ws4.Range("I21").Value = 0.118
$ws4.Range("J21").Value = 0.001
$ws4.Range("K21").Value = 0.02

$ws4.Range("B22").Value = 3.573
$ws4.Range("C22").Value = 2.1
$ws4.Range("D22").Value = 5.445
$ws4.Range("E22").Value = 0.044
$ws4.Range("F22").Value = 0.0094
$ws4.Range("G22").Value = 0.5479000000000001
$ws4.Range("H22").Value = 0.44
$ws4.Range("I22").Value = 0.102

$ws4.Range("B23").Value = 0.057
$ws4.Range("D23").Value = 0.197
$ws4.Range("E23").Value = 0.0271
$ws4.Range("F23").Value = 0.0011
$ws4.Range("G23").Value = 0.5395
$ws4.Range("H23").Value = 0.77
$ws4.Range("I23").Value = 0.26
$ws4.Range("J23").Value = 0.453

$ws4.Range("B24").Value = 0.021
$ws4.Range("C24").Value = 0.001
$ws4.Range("D24").Value = 0.077
$ws4.Range("E24").Value = 0.011
$ws4.Range("F24").Value = 0
$ws4.Range("G24").Value = 0.0461

$ws4.Range("B25").Value = 3.573
$ws4.Range("C25").Value = 2.1
$ws4.Range("D25").Value = 5.445
$ws4.Range("F25").Value = 0.0118
$ws4.Range("G25").Value = 8.953900000000001
$ws4.Range("H25").Value = 1.76
$ws4.Range("I25").Value = 0.5659999999999999
$ws4.Range("K25").Value = 21.69
